$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / date text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Row 14 ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J14").Value = 3

# --- Row 15 ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 155.555555555556
$ws.Range("L15").Value = 155.555555555556

# --- Row 16 ---
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -77.777777777777
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = -30.15873015873
$ws.Range("L16").Value = -27.868852459016

# --- Row 17 ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -32
$ws.Range("I17").Value = 187
$ws.Range("J17").Value = 167
$ws.Range("K17").Value = 11.976047904191
$ws.Range("L17").Value = 6.857142857142

# --- Row 18 ---
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("E18").Value = 100
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -24.528301886792
$ws.Range("L18").Value = -27.272727272727

# --- Row 19 ---
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -11.764705882352
$ws.Range("I19").Value = 154
$ws.Range("J19").Value = 189
$ws.Range("K19").Value = -18.518518518518
$ws.Range("L19").Value = -30.316742081448

# --- Row 20 ---
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -68.421052631578
$ws.Range("I20").Value = 118
$ws.Range("J20").Value = 125
$ws.Range("K20").Value = -5.6
$ws.Range("L20").Value = -19.17808219178

# --- Row 21 ---
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -52.380952380952
$ws.Range("F21").Value = 47
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -38.157894736842
$ws.Range("I21").Value = 566
$ws.Range("J21").Value = 609
$ws.Range("K21").Value = -7.060755336617
$ws.Range("L21").Value = -15.396113602391

# --- Row 24 ---
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = 113.157894736842
$ws.Range("I24").Value = 418
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 4.5
$ws.Range("L24").Value = -3.24074074074

# --- Row 25 ---
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 73
$ws.Range("J25").Value = 112
$ws.Range("K25").Value = -34.821428571428
$ws.Range("L25").Value = -38.135593220339

# --- Row 26 ---
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 19.354838709677
$ws.Range("I26").Value = 333
$ws.Range("J26").Value = 318
$ws.Range("K26").Value = 4.716981132075
$ws.Range("L26").Value = 30.588235294117

# --- Row 27 ---
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 64.705882352941
$ws.Range("L27").Value = 75

# --- Row 29 ---
$ws.Range("D29").Value = 5
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G29").Value = 5
$ws.Range("G29").NumberFormat = '#,##0'
$ws.Range("H29").Value = -80
$ws.Range("H29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -84.615384615384

# --- Row 30 ---
$ws.Range("D30").Value = 2
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G30").Value = 2
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("H30").Value = -50
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -75

# --- Cells reverting to "no data" text placeholders (copy style+value from a stable source cell) ---
$ws.Range("C22").Copy()
$ws.Range("C29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C22").Copy()
$ws.Range("C29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("C22").Copy()
$ws.Range("C30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C22").Copy()
$ws.Range("C30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("C22").Copy()
$ws.Range("F33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C22").Copy()
$ws.Range("F33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$excel.CutCopyMode = $false
